$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns D and E keep their text formatting so values like "1.000"
# are not reinterpreted as numbers when assigned.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '29.134.85'
$ws.Range('E2').Value = '  -0.35%  '
$ws.Range('D3').Value = '1.853.59'
$ws.Range('E3').Value = '  +0.09%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '238.14'
$ws.Range('E5').Value = '  -0.22%  '
$ws.Range('D6').Value = '0.6895'
$ws.Range('E6').Value = '  -0.90%  '
$ws.Range('D7').Value = '1.000'
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = '0.07728'
$ws.Range('E8').Value = '  +1.54%  '
$ws.Range('D9').Value = '0.3035'
$ws.Range('E9').Value = '  -1.29%  '
$ws.Range('D10').Value = '23.13'
$ws.Range('E10').Value = '  -1.80%  '
$ws.Range('D11').Value = '0.08077'
$ws.Range('E11').Value = '  -0.18%  '
$ws.Range('D12').Value = '1.885.80'
$ws.Range('E12').Value = '  +3.41%  '
$ws.Range('D13').Value = '0.7209'
$ws.Range('E13').Value = '  -0.34%  '
$ws.Range('D14').Value = '5.191'
$ws.Range('E14').Value = '  -0.21%  '
$ws.Range('D15').Value = '89.37'
$ws.Range('E15').Value = '  +0.11%  '
$ws.Range('D16').Value = '29.140.35'
$ws.Range('E16').Value = '  -0.12%  '
$ws.Range('D17').Value = '5.731'
$ws.Range('E17').Value = '  -2.50%  '
$ws.Range('B18').Value = 'Avalanche'
$ws.Range('C18').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D18').Value = '13.25'
$ws.Range('E18').Value = '  +1.16%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').Value = '0.000007775'
$ws.Range('E19').Value = '  +0.45%  '
$ws.Range('D20').Value = '234.28'
$ws.Range('E20').Value = '  -3.33%  '
$ws.Range('D21').Value = '1.000'
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('D22').Value = '2.107.88'
$ws.Range('E22').Value = '  +1.34%  '
$ws.Range('D23').Value = '1.001'
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('D24').Value = '7.473'
$ws.Range('E24').Value = '  -1.80%  '
$ws.Range('D25').Value = '161.55'
$ws.Range('E25').Value = '  -1.02%  '
$ws.Range('D26').Value = '8.982'
$ws.Range('E26').Value = '  -0.95%  '
$ws.Range('D27').Value = '0.1427'
$ws.Range('E27').Value = '  -2.18%  '
$ws.Range('E28').Value = '  -0.38%  '
$ws.Range('D29').Value = '1.952'
$ws.Range('E29').Value = '  +0.84%  '
$ws.Range('D30').Value = '1.413'
$ws.Range('E30').Value = '  +0.95%  '
$ws.Range('D31').Value = '4.495'
$ws.Range('E31').Value = '  +1.39%  '
$ws.Range('D32').Value = '1.482'
$ws.Range('E32').Value = '  -1.46%  '
$ws.Range('E33').Value = '  -0.67%  '
$ws.Range('D34').Value = '0.05200'
$ws.Range('E34').Value = '  -1.48%  '
$ws.Range('E35').Value = '  -1.51%  '
$ws.Range('D36').Value = '0.7057'
$ws.Range('E36').Value = '  -1.05%  '
$ws.Range('D37').Value = '0.9990'
$ws.Range('E37').Value = '  -0.36%  '
$ws.Range('D38').Value = '2.659'
$ws.Range('E38').Value = '  -0.17%  '
$ws.Range('D39').Value = '0.01851'
$ws.Range('E39').Value = '  -0.48%  '
$ws.Range('D40').Value = '2.715'
$ws.Range('E40').Value = '  +0.61%  '
$ws.Range('D41').Value = '0.9288'
$ws.Range('E41').Value = '  -0.74%  '
$ws.Range('D42').Value = '1.104.24'
$ws.Range('E42').Value = '  +5.75%  '
$ws.Range('D43').Value = '0.4283'
$ws.Range('E43').Value = '  -0.41%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').Value = '70.53'
$ws.Range('E44').Value = '  +1.72%  '
$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D45').Value = '5.866'
$ws.Range('E45').Value = '  -0.32%  '
$ws.Range('D46').Value = '1.0000'
$ws.Range('E46').Value = '  -0.11%  '
$ws.Range('D47').Value = '102.75'
$ws.Range('E47').Value = '  +0.24%  '
$ws.Range('D48').Value = '1.800'
$ws.Range('E48').Value = '  +3.38%  '
$ws.Range('D49').Value = '2.006.61'
$ws.Range('E49').Value = '  +1.69%  '
$ws.Range('D50').Value = '9.141'
$ws.Range('E50').Value = '  -1.45%  '
$ws.Range('D51').Value = '6.991'
$ws.Range('E51').Value = '  -3.72%  '
